$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header D1: cargo -> UserRole
$ws.Cells.Item(1, 4).Value = "UserRole"

# Adjust row heights (rows resized slightly due to content width changes)
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5

# Update trilha (C) and UserRole (D) values for each user row
$ws.Cells.Item(2, 3).Value = "BUSINESSMAN"
$ws.Cells.Item(2, 4).Value = "EMPLOYER"
$ws.Cells.Item(3, 3).Value = "BUSINESSMAN"
$ws.Cells.Item(3, 4).Value = "LEADER"
$ws.Cells.Item(4, 3).Value = "DESIGNER"
$ws.Cells.Item(4, 4).Value = "EMPLOYER"
$ws.Cells.Item(5, 3).Value = "BUSINESSMAN"
$ws.Cells.Item(5, 4).Value = "MANAGER"
$ws.Cells.Item(6, 3).Value = "DESIGNER"
$ws.Cells.Item(6, 4).Value = "EMPLOYER"
$ws.Cells.Item(7, 3).Value = "DEVELOPER"
$ws.Cells.Item(7, 4).Value = "MENTOR"
$ws.Cells.Item(8, 3).Value = "ARCHITECT"
$ws.Cells.Item(8, 4).Value = "MENTOR"
$ws.Cells.Item(9, 3).Value = "DEVELOPER"
$ws.Cells.Item(9, 4).Value = "EMPLOYER"
$ws.Cells.Item(10, 3).Value = "LEADER"
$ws.Cells.Item(10, 4).Value = "LEADER"
$ws.Cells.Item(11, 3).Value = "LEADER"
$ws.Cells.Item(11, 4).Value = "LEADER"
$ws.Cells.Item(12, 3).Value = "DEVELOPER"
$ws.Cells.Item(12, 4).Value = "LEADER"
$ws.Cells.Item(13, 3).Value = "BUSINESSMAN"
$ws.Cells.Item(13, 4).Value = "EMPLOYER"
$ws.Cells.Item(14, 3).Value = "MANAGER"
$ws.Cells.Item(14, 4).Value = "LEADER"
$ws.Cells.Item(15, 3).Value = "DESIGNER"
$ws.Cells.Item(15, 4).Value = "EMPLOYER"
$ws.Cells.Item(16, 3).Value = "DEVELOPER"
$ws.Cells.Item(16, 4).Value = "EMPLOYER"
$ws.Cells.Item(17, 3).Value = "DEVELOPER"
$ws.Cells.Item(17, 4).Value = "EMPLOYER"
$ws.Cells.Item(18, 3).Value = "DESIGNER"
$ws.Cells.Item(18, 4).Value = "EMPLOYER"
$ws.Cells.Item(19, 3).Value = "DESIGNER"
$ws.Cells.Item(19, 4).Value = "LEADER"
$ws.Cells.Item(20, 3).Value = "DEVELOPER"
$ws.Cells.Item(20, 4).Value = "EMPLOYER"
$ws.Cells.Item(21, 3).Value = "DEVELOPER"
$ws.Cells.Item(21, 4).Value = "LEADER"
$ws.Cells.Item(22, 3).Value = "MANAGER"
$ws.Cells.Item(22, 4).Value = "LEADER"
$ws.Cells.Item(23, 3).Value = "DEVELOPER"
$ws.Cells.Item(23, 4).Value = "EMPLOYER"
$ws.Cells.Item(24, 3).Value = "ARCHITECT"
$ws.Cells.Item(24, 4).Value = "LEADER"
$ws.Cells.Item(25, 3).Value = "DEVELOPER"
$ws.Cells.Item(25, 4).Value = "LEADER"
$ws.Cells.Item(26, 3).Value = "DEVELOPER"
$ws.Cells.Item(26, 4).Value = "EMPLOYER"
$ws.Cells.Item(27, 3).Value = "LEADER"
$ws.Cells.Item(27, 4).Value = "LEADER"
$ws.Cells.Item(28, 3).Value = "MANAGER"
$ws.Cells.Item(28, 4).Value = "EMPLOYER"
$ws.Cells.Item(29, 3).Value = "BUSINESSMAN"
$ws.Cells.Item(29, 4).Value = "LEADER"
$ws.Cells.Item(30, 3).Value = "LEADER"
$ws.Cells.Item(30, 4).Value = "LEADER"
$ws.Cells.Item(31, 3).Value = "DEVELOPER"
$ws.Cells.Item(31, 4).Value = "EMPLOYER"
$ws.Cells.Item(32, 3).Value = "BUSINESSMAN"
$ws.Cells.Item(32, 4).Value = "EMPLOYER"
$ws.Cells.Item(33, 3).Value = "DEVELOPER"
$ws.Cells.Item(33, 4).Value = "EMPLOYER"
$ws.Cells.Item(34, 3).Value = "DEVELOPER"
$ws.Cells.Item(34, 4).Value = "LEADER"
$ws.Cells.Item(35, 3).Value = "DEVELOPER"
$ws.Cells.Item(35, 4).Value = "EMPLOYER"
$ws.Cells.Item(36, 3).Value = "DEVELOPER"
$ws.Cells.Item(36, 4).Value = "EMPLOYER"
$ws.Cells.Item(37, 3).Value = "LEADER"
$ws.Cells.Item(37, 4).Value = "LEADER"
$ws.Cells.Item(38, 3).Value = "BUSINESSMAN"
$ws.Cells.Item(38, 4).Value = "EMPLOYER"
$ws.Cells.Item(39, 3).Value = "MANAGER"
$ws.Cells.Item(39, 4).Value = "LEADER"
